$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update existing "GDP per Capita" values (years 1800-2010, column E) ---
# These values are stored as text (not numbers) in the source data, matching the
# original workbook convention, so NumberFormat is set to text ("@") before writing.
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "1484"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "1591"
$ws.Cells.Item(52, 5).NumberFormat = "@"
$ws.Cells.Item(52, 5).Value = "1994"
$ws.Cells.Item(62, 5).NumberFormat = "@"
$ws.Cells.Item(62, 5).Value = "2160"
$ws.Cells.Item(72, 5).NumberFormat = "@"
$ws.Cells.Item(72, 5).Value = "2340"
$ws.Cells.Item(77, 5).NumberFormat = "@"
$ws.Cells.Item(77, 5).Value = "2606"
$ws.Cells.Item(78, 5).NumberFormat = "@"
$ws.Cells.Item(78, 5).Value = "2611"
$ws.Cells.Item(79, 5).NumberFormat = "@"
$ws.Cells.Item(79, 5).Value = "2820"
$ws.Cells.Item(80, 5).NumberFormat = "@"
$ws.Cells.Item(80, 5).Value = "2619"
$ws.Cells.Item(81, 5).NumberFormat = "@"
$ws.Cells.Item(81, 5).Value = "2664"
$ws.Cells.Item(82, 5).NumberFormat = "@"
$ws.Cells.Item(82, 5).Value = "2557"
$ws.Cells.Item(83, 5).NumberFormat = "@"
$ws.Cells.Item(83, 5).Value = "2541"
$ws.Cells.Item(84, 5).NumberFormat = "@"
$ws.Cells.Item(84, 5).Value = "3108"
$ws.Cells.Item(85, 5).NumberFormat = "@"
$ws.Cells.Item(85, 5).Value = "3373"
$ws.Cells.Item(86, 5).NumberFormat = "@"
$ws.Cells.Item(86, 5).Value = "3486"
$ws.Cells.Item(87, 5).NumberFormat = "@"
$ws.Cells.Item(87, 5).Value = "3904"
$ws.Cells.Item(88, 5).NumberFormat = "@"
$ws.Cells.Item(88, 5).Value = "3771"
$ws.Cells.Item(89, 5).NumberFormat = "@"
$ws.Cells.Item(89, 5).Value = "3841"
$ws.Cells.Item(90, 5).NumberFormat = "@"
$ws.Cells.Item(90, 5).Value = "4219"
$ws.Cells.Item(91, 5).NumberFormat = "@"
$ws.Cells.Item(91, 5).Value = "4288"
$ws.Cells.Item(92, 5).NumberFormat = "@"
$ws.Cells.Item(92, 5).Value = "3851"
$ws.Cells.Item(93, 5).NumberFormat = "@"
$ws.Cells.Item(93, 5).Value = "3631"
$ws.Cells.Item(94, 5).NumberFormat = "@"
$ws.Cells.Item(94, 5).Value = "4251"
$ws.Cells.Item(95, 5).NumberFormat = "@"
$ws.Cells.Item(95, 5).Value = "4409"
$ws.Cells.Item(96, 5).NumberFormat = "@"
$ws.Cells.Item(96, 5).Value = "4968"
$ws.Cells.Item(97, 5).NumberFormat = "@"
$ws.Cells.Item(97, 5).Value = "5384"
$ws.Cells.Item(98, 5).NumberFormat = "@"
$ws.Cells.Item(98, 5).Value = "5716"
$ws.Cells.Item(99, 5).NumberFormat = "@"
$ws.Cells.Item(99, 5).Value = "4500"
$ws.Cells.Item(100, 5).NumberFormat = "@"
$ws.Cells.Item(100, 5).Value = "4741"
$ws.Cells.Item(101, 5).NumberFormat = "@"
$ws.Cells.Item(101, 5).Value = "5412"
$ws.Cells.Item(102, 5).NumberFormat = "@"
$ws.Cells.Item(102, 5).Value = "4583"
$ws.Cells.Item(103, 5).NumberFormat = "@"
$ws.Cells.Item(103, 5).Value = "4591"
$ws.Cells.Item(104, 5).NumberFormat = "@"
$ws.Cells.Item(104, 5).Value = "4331"
$ws.Cells.Item(105, 5).NumberFormat = "@"
$ws.Cells.Item(105, 5).Value = "4769"
$ws.Cells.Item(106, 5).NumberFormat = "@"
$ws.Cells.Item(106, 5).Value = "5086"
$ws.Cells.Item(107, 5).NumberFormat = "@"
$ws.Cells.Item(107, 5).Value = "5545"
$ws.Cells.Item(108, 5).NumberFormat = "@"
$ws.Cells.Item(108, 5).Value = "5608"
$ws.Cells.Item(109, 5).NumberFormat = "@"
$ws.Cells.Item(109, 5).Value = "5514"
$ws.Cells.Item(110, 5).NumberFormat = "@"
$ws.Cells.Item(110, 5).Value = "5829"
$ws.Cells.Item(111, 5).NumberFormat = "@"
$ws.Cells.Item(111, 5).Value = "5896"
$ws.Cells.Item(112, 5).NumberFormat = "@"
$ws.Cells.Item(112, 5).Value = "6092"
$ws.Cells.Item(113, 5).NumberFormat = "@"
$ws.Cells.Item(113, 5).Value = "5971"
$ws.Cells.Item(114, 5).NumberFormat = "@"
$ws.Cells.Item(114, 5).Value = "6223"
$ws.Cells.Item(115, 5).NumberFormat = "@"
$ws.Cells.Item(115, 5).Value = "6052"
$ws.Cells.Item(116, 5).NumberFormat = "@"
$ws.Cells.Item(116, 5).Value = "5263"
$ws.Cells.Item(117, 5).NumberFormat = "@"
$ws.Cells.Item(117, 5).Value = "5171"
$ws.Cells.Item(118, 5).NumberFormat = "@"
$ws.Cells.Item(118, 5).Value = "4927"
$ws.Cells.Item(119, 5).NumberFormat = "@"
$ws.Cells.Item(119, 5).Value = "4447"
$ws.Cells.Item(120, 5).NumberFormat = "@"
$ws.Cells.Item(120, 5).Value = "5177"
$ws.Cells.Item(121, 5).NumberFormat = "@"
$ws.Cells.Item(121, 5).Value = "5271"
$ws.Cells.Item(122, 5).NumberFormat = "@"
$ws.Cells.Item(122, 5).Value = "5536"
$ws.Cells.Item(123, 5).NumberFormat = "@"
$ws.Cells.Item(123, 5).Value = "5533"
$ws.Cells.Item(124, 5).NumberFormat = "@"
$ws.Cells.Item(124, 5).Value = "5796"
$ws.Cells.Item(125, 5).NumberFormat = "@"
$ws.Cells.Item(125, 5).Value = "6213"
$ws.Cells.Item(126, 5).NumberFormat = "@"
$ws.Cells.Item(126, 5).Value = "6464"
$ws.Cells.Item(127, 5).NumberFormat = "@"
$ws.Cells.Item(127, 5).Value = "6247"
$ws.Cells.Item(128, 5).NumberFormat = "@"
$ws.Cells.Item(128, 5).Value = "6366"
$ws.Cells.Item(129, 5).NumberFormat = "@"
$ws.Cells.Item(129, 5).Value = "6625"
$ws.Cells.Item(130, 5).NumberFormat = "@"
$ws.Cells.Item(130, 5).Value = "6840"
$ws.Cells.Item(131, 5).NumberFormat = "@"
$ws.Cells.Item(131, 5).Value = "6961"
$ws.Cells.Item(132, 5).NumberFormat = "@"
$ws.Cells.Item(132, 5).Value = "6503"
$ws.Cells.Item(133, 5).NumberFormat = "@"
$ws.Cells.Item(133, 5).Value = "5917"
$ws.Cells.Item(134, 5).NumberFormat = "@"
$ws.Cells.Item(134, 5).Value = "5614"
$ws.Cells.Item(135, 5).NumberFormat = "@"
$ws.Cells.Item(135, 5).Value = "5772"
$ws.Cells.Item(136, 5).NumberFormat = "@"
$ws.Cells.Item(136, 5).Value = "6129"
$ws.Cells.Item(137, 5).NumberFormat = "@"
$ws.Cells.Item(137, 5).Value = "6296"
$ws.Cells.Item(138, 5).NumberFormat = "@"
$ws.Cells.Item(138, 5).Value = "6236"
$ws.Cells.Item(139, 5).NumberFormat = "@"
$ws.Cells.Item(139, 5).Value = "6575"
$ws.Cells.Item(140, 5).NumberFormat = "@"
$ws.Cells.Item(140, 5).Value = "6491"
$ws.Cells.Item(141, 5).NumberFormat = "@"
$ws.Cells.Item(141, 5).Value = "6612"
$ws.Cells.Item(142, 5).NumberFormat = "@"
$ws.Cells.Item(142, 5).Value = "6633"
$ws.Cells.Item(143, 5).NumberFormat = "@"
$ws.Cells.Item(143, 5).Value = "6861"
$ws.Cells.Item(144, 5).NumberFormat = "@"
$ws.Cells.Item(144, 5).Value = "6829"
$ws.Cells.Item(145, 5).NumberFormat = "@"
$ws.Cells.Item(145, 5).Value = "6666"
$ws.Cells.Item(146, 5).NumberFormat = "@"
$ws.Cells.Item(146, 5).Value = "7299"
$ws.Cells.Item(147, 5).NumberFormat = "@"
$ws.Cells.Item(147, 5).Value = "6943"
$ws.Cells.Item(148, 5).NumberFormat = "@"
$ws.Cells.Item(148, 5).Value = "7436"
$ws.Cells.Item(149, 5).NumberFormat = "@"
$ws.Cells.Item(149, 5).Value = "8112"
$ws.Cells.Item(150, 5).NumberFormat = "@"
$ws.Cells.Item(150, 5).Value = "8372"
$ws.Cells.Item(151, 5).NumberFormat = "@"
$ws.Cells.Item(151, 5).Value = "8045"
$ws.Cells.Item(152, 5).NumberFormat = "@"
$ws.Cells.Item(152, 5).Value = "7949"
$ws.Cells.Item(153, 5).NumberFormat = "@"
$ws.Cells.Item(153, 5).Value = "8086"
$ws.Cells.Item(154, 5).NumberFormat = "@"
$ws.Cells.Item(154, 5).Value = "7519"
$ws.Cells.Item(155, 5).NumberFormat = "@"
$ws.Cells.Item(155, 5).Value = "7769"
$ws.Cells.Item(156, 5).NumberFormat = "@"
$ws.Cells.Item(156, 5).Value = "7938"
$ws.Cells.Item(157, 5).NumberFormat = "@"
$ws.Cells.Item(157, 5).Value = "8348"
$ws.Cells.Item(158, 5).NumberFormat = "@"
$ws.Cells.Item(158, 5).Value = "8424"
$ws.Cells.Item(159, 5).NumberFormat = "@"
$ws.Cells.Item(159, 5).Value = "8705"
$ws.Cells.Item(160, 5).NumberFormat = "@"
$ws.Cells.Item(160, 5).Value = "9083"
$ws.Cells.Item(161, 5).NumberFormat = "@"
$ws.Cells.Item(161, 5).Value = "8354"
$ws.Cells.Item(162, 5).NumberFormat = "@"
$ws.Cells.Item(162, 5).Value = "8861"
$ws.Cells.Item(163, 5).NumberFormat = "@"
$ws.Cells.Item(163, 5).Value = "9344"
$ws.Cells.Item(164, 5).NumberFormat = "@"
$ws.Cells.Item(164, 5).Value = "9049"
$ws.Cells.Item(165, 5).NumberFormat = "@"
$ws.Cells.Item(165, 5).Value = "8695"
$ws.Cells.Item(166, 5).NumberFormat = "@"
$ws.Cells.Item(166, 5).Value = "9446"
$ws.Cells.Item(167, 5).NumberFormat = "@"
$ws.Cells.Item(167, 5).Value = "10155"
$ws.Cells.Item(168, 5).NumberFormat = "@"
$ws.Cells.Item(168, 5).Value = "10076"
$ws.Cells.Item(169, 5).NumberFormat = "@"
$ws.Cells.Item(169, 5).Value = "10200"
$ws.Cells.Item(170, 5).NumberFormat = "@"
$ws.Cells.Item(170, 5).Value = "10485"
$ws.Cells.Item(171, 5).NumberFormat = "@"
$ws.Cells.Item(171, 5).Value = "11217"
$ws.Cells.Item(172, 5).NumberFormat = "@"
$ws.Cells.Item(172, 5).Value = "11639"
$ws.Cells.Item(173, 5).NumberFormat = "@"
$ws.Cells.Item(173, 5).Value = "12003"
$ws.Cells.Item(174, 5).NumberFormat = "@"
$ws.Cells.Item(174, 5).Value = "12170"
$ws.Cells.Item(175, 5).NumberFormat = "@"
$ws.Cells.Item(175, 5).Value = "12691"
$ws.Cells.Item(176, 5).NumberFormat = "@"
$ws.Cells.Item(176, 5).Value = "13284"
$ws.Cells.Item(177, 5).NumberFormat = "@"
$ws.Cells.Item(177, 5).Value = "12946"
$ws.Cells.Item(178, 5).NumberFormat = "@"
$ws.Cells.Item(178, 5).Value = "12696"
$ws.Cells.Item(179, 5).NumberFormat = "@"
$ws.Cells.Item(179, 5).Value = "13236"
$ws.Cells.Item(180, 5).NumberFormat = "@"
$ws.Cells.Item(180, 5).Value = "12444"
$ws.Cells.Item(181, 5).NumberFormat = "@"
$ws.Cells.Item(181, 5).Value = "13114"
$ws.Cells.Item(182, 5).NumberFormat = "@"
$ws.Cells.Item(182, 5).Value = "13080"
$ws.Cells.Item(183, 5).NumberFormat = "@"
$ws.Cells.Item(183, 5).Value = "12125"
$ws.Cells.Item(184, 5).NumberFormat = "@"
$ws.Cells.Item(184, 5).Value = "11550"
$ws.Cells.Item(185, 5).NumberFormat = "@"
$ws.Cells.Item(185, 5).Value = "11775"
$ws.Cells.Item(186, 5).NumberFormat = "@"
$ws.Cells.Item(186, 5).Value = "11837"
$ws.Cells.Item(187, 5).NumberFormat = "@"
$ws.Cells.Item(187, 5).Value = "10895"
$ws.Cells.Item(188, 5).NumberFormat = "@"
$ws.Cells.Item(188, 5).Value = "11515"
$ws.Cells.Item(189, 5).NumberFormat = "@"
$ws.Cells.Item(189, 5).Value = "11633"
$ws.Cells.Item(190, 5).NumberFormat = "@"
$ws.Cells.Item(190, 5).Value = "11244"
$ws.Cells.Item(191, 5).NumberFormat = "@"
$ws.Cells.Item(191, 5).Value = "10393"
$ws.Cells.Item(192, 5).NumberFormat = "@"
$ws.Cells.Item(192, 5).Value = "10254"
$ws.Cells.Item(193, 5).NumberFormat = "@"
$ws.Cells.Item(193, 5).Value = "11223.9846365402"
$ws.Cells.Item(194, 5).NumberFormat = "@"
$ws.Cells.Item(194, 5).Value = "12267.0835867812"
$ws.Cells.Item(195, 5).NumberFormat = "@"
$ws.Cells.Item(195, 5).Value = "12926.7946652176"
$ws.Cells.Item(196, 5).NumberFormat = "@"
$ws.Cells.Item(196, 5).Value = "13571.379431201"
$ws.Cells.Item(197, 5).NumberFormat = "@"
$ws.Cells.Item(197, 5).Value = "13086.0364515543"
$ws.Cells.Item(198, 5).NumberFormat = "@"
$ws.Cells.Item(198, 5).Value = "13715.2268187077"
$ws.Cells.Item(199, 5).NumberFormat = "@"
$ws.Cells.Item(199, 5).Value = "14722.6459400365"
$ws.Cells.Item(200, 5).NumberFormat = "@"
$ws.Cells.Item(200, 5).Value = "15185.9586884861"
$ws.Cells.Item(201, 5).NumberFormat = "@"
$ws.Cells.Item(201, 5).Value = "14577.8322730293"
$ws.Cells.Item(202, 5).NumberFormat = "@"
$ws.Cells.Item(202, 5).Value = "14368.9427373623"
$ws.Cells.Item(203, 5).NumberFormat = "@"
$ws.Cells.Item(203, 5).Value = "13651.980987935"
$ws.Cells.Item(204, 5).NumberFormat = "@"
$ws.Cells.Item(204, 5).Value = "12094.7728377935"
$ws.Cells.Item(205, 5).NumberFormat = "@"
$ws.Cells.Item(205, 5).Value = "13088.5650108326"
$ws.Cells.Item(206, 5).NumberFormat = "@"
$ws.Cells.Item(206, 5).Value = "14183.225137777"
$ws.Cells.Item(207, 5).NumberFormat = "@"
$ws.Cells.Item(207, 5).Value = "15344.1571340567"
$ws.Cells.Item(208, 5).NumberFormat = "@"
$ws.Cells.Item(208, 5).Value = "16490.4735815761"
$ws.Cells.Item(209, 5).NumberFormat = "@"
$ws.Cells.Item(209, 5).Value = "17891.9434108489"
$ws.Cells.Item(210, 5).NumberFormat = "@"
$ws.Cells.Item(210, 5).Value = "18520.3040835795"
$ws.Cells.Item(211, 5).NumberFormat = "@"
$ws.Cells.Item(211, 5).Value = "17328.97575242"
$ws.Cells.Item(212, 5).NumberFormat = "@"
$ws.Cells.Item(212, 5).Value = "18979.9917029921"

# --- Add new rows for years 2011-2016 ---
$ws.Cells.Item(213, 1).Value = 32
$ws.Cells.Item(213, 2).Value = "Argentina"
$ws.Cells.Item(213, 3).Value = "GDP per Capita"
$ws.Cells.Item(213, 4).Value = 2011
$ws.Cells.Item(213, 5).NumberFormat = "@"
$ws.Cells.Item(213, 5).Value = "20003"
$ws.Cells.Item(214, 1).Value = 32
$ws.Cells.Item(214, 2).Value = "Argentina"
$ws.Cells.Item(214, 3).Value = "GDP per Capita"
$ws.Cells.Item(214, 4).Value = 2012
$ws.Cells.Item(214, 5).NumberFormat = "@"
$ws.Cells.Item(214, 5).Value = "19599"
$ws.Cells.Item(215, 1).Value = 32
$ws.Cells.Item(215, 2).Value = "Argentina"
$ws.Cells.Item(215, 3).Value = "GDP per Capita"
$ws.Cells.Item(215, 4).Value = 2013
$ws.Cells.Item(215, 5).NumberFormat = "@"
$ws.Cells.Item(215, 5).Value = "19873"
$ws.Cells.Item(216, 1).Value = 32
$ws.Cells.Item(216, 2).Value = "Argentina"
$ws.Cells.Item(216, 3).Value = "GDP per Capita"
$ws.Cells.Item(216, 4).Value = 2014
$ws.Cells.Item(216, 5).NumberFormat = "@"
$ws.Cells.Item(216, 5).Value = "19183"
$ws.Cells.Item(217, 1).Value = 32
$ws.Cells.Item(217, 2).Value = "Argentina"
$ws.Cells.Item(217, 3).Value = "GDP per Capita"
$ws.Cells.Item(217, 4).Value = 2015
$ws.Cells.Item(217, 5).NumberFormat = "@"
$ws.Cells.Item(217, 5).Value = "19502"
$ws.Cells.Item(218, 1).Value = 32
$ws.Cells.Item(218, 2).Value = "Argentina"
$ws.Cells.Item(218, 3).Value = "GDP per Capita"
$ws.Cells.Item(218, 4).Value = 2016
$ws.Cells.Item(218, 5).NumberFormat = "@"
$ws.Cells.Item(218, 5).Value = "18875"